$d = $word.ActiveDocument

# Update the header date paragraph.
$d.Content.Find.Execute("2025-09-28 Sunday", $true, $false, $false, $false,
                         $false, $true, 1, $false, "2025-09-29 Monday", 2)

$t = $d.Tables.Item(1)

# Row 1 (division problems)
$t.Cell(1, 1).Range.Text = "19÷6=3, 1"
$t.Cell(1, 2).Range.Text = "91÷2=45, 1"
$t.Cell(1, 3).Range.Text = "36÷3=12, 0"
$t.Cell(1, 4).Range.Text = "57÷6=9, 3"
$t.Cell(1, 5).Range.Text = "50÷2=25, 0"

# Row 5
$t.Cell(5, 1).Range.Text = "89÷3=29, 2"
$t.Cell(5, 2).Range.Text = "32÷3=10, 2"
$t.Cell(5, 3).Range.Text = "41÷6=6, 5"
$t.Cell(5, 4).Range.Text = "64÷6=10, 4"
$t.Cell(5, 5).Range.Text = "51÷7=7, 2"

# Row 9
$t.Cell(9, 1).Range.Text = "38÷7=5, 3"
$t.Cell(9, 2).Range.Text = "29÷3=9, 2"
$t.Cell(9, 3).Range.Text = "53÷6=8, 5"
$t.Cell(9, 4).Range.Text = "50÷2=25, 0"
$t.Cell(9, 5).Range.Text = "29÷4=7, 1"

# Row 13
$t.Cell(13, 1).Range.Text = "71÷5=14, 1"
$t.Cell(13, 2).Range.Text = "47÷6=7, 5"
$t.Cell(13, 3).Range.Text = "98÷4=24, 2"
$t.Cell(13, 4).Range.Text = "85÷2=42, 1"
$t.Cell(13, 5).Range.Text = "68÷9=7, 5"

# Row 17
$t.Cell(17, 1).Range.Text = "74÷6=12, 2"
$t.Cell(17, 2).Range.Text = "89÷2=44, 1"
$t.Cell(17, 3).Range.Text = "96÷4=24, 0"
$t.Cell(17, 4).Range.Text = "80÷4=20, 0"
$t.Cell(17, 5).Range.Text = "16÷5=3, 1"
